# Fruta / hortaliza, semanal
#
# A new weekly price-report row is inserted at row 209 (pushing the
# existing rows 209:278 down to 210:279). The new row carries a fresh
# reading (fecha 2022-01-27 / serial 44588) while everything that used
# to live at 209:278 simply slides down one row, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 209; Excel shifts 209:278 -> 210:279
# and, by default, the new row inherits formatting from the row above it
# (which is what keeps column D's date style intact).
$ws.Rows(209).Insert()

# Most fields of the new record are identical to the (now-shifted) row
# directly below it, so copy those across first ...
$ws.Range("A209").Value = $ws.Range("A210").Value()
$ws.Range("B209").Value = $ws.Range("B210").Value()
$ws.Range("C209").Value = $ws.Range("C210").Value()
$ws.Range("E209").Value = $ws.Range("E210").Value()
$ws.Range("F209").Value = $ws.Range("F210").Value()
$ws.Range("G209").Value = $ws.Range("G210").Value()
$ws.Range("H209").Value = $ws.Range("H210").Value()
$ws.Range("I209").Value = $ws.Range("I210").Value()
$ws.Range("N209").Value = $ws.Range("N210").Value()
$ws.Range("Q209").Value = $ws.Range("Q210").Value()
$ws.Range("R209").Value = $ws.Range("R210").Value()

# ... then fill in the values that actually differ for this new record.
$ws.Range("D209").Value = 44588
$ws.Range("J209").Value = 300
$ws.Range("K209").Value = 8000
$ws.Range("L209").Value = 8000
$ws.Range("M209").Value = 8000
$ws.Range("O209").Value = "Región del Maule"
$ws.Range("P209").Value = 160
